$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill Summary")

function Set-TextCell($sheet, $addr, $text) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
}

# Row 8
Set-TextCell $ws 'A8' ''
$ws.Range("C8").Value = 48
Set-TextCell $ws 'D8' '1.0'
Set-TextCell $ws 'E8' 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
Set-TextCell $ws 'G8' '0.00'

# Row 9
$ws.Range("C9").Value = 4
Set-TextCell $ws 'D9' '2'
Set-TextCell $ws 'E9' 'Short point (up to 3 mtr.)'
$ws.Range("F9").Value = 256
Set-TextCell $ws 'G9' '1024.00'

# Row 10
$ws.Range("C10").Value = 15
Set-TextCell $ws 'G10' '9930.00'

# Row 11
$ws.Range("C11").Value = 19

# Row 12
$ws.Range("C12").Value = 25
Set-TextCell $ws 'G12' '3400.00'

# Row 13
$ws.Range("C13").Value = 28
Set-TextCell $ws 'D13' '4.0'
Set-TextCell $ws 'E13' 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 50
Set-TextCell $ws 'G13' '1400.00'

# Row 14
$ws.Range("C14").Value = 33
Set-TextCell $ws 'D14' '6.0'
Set-TextCell $ws 'E14' 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 78
Set-TextCell $ws 'G14' '2574.00'

# Row 15
Set-TextCell $ws 'D15' '7.0'
Set-TextCell $ws 'E15' 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 30
Set-TextCell $ws 'G15' '690.00'

# Row 16
$ws.Range("C16").Value = 49
Set-TextCell $ws 'D16' '8.0'
Set-TextCell $ws 'E16' 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 30
Set-TextCell $ws 'G16' '1470.00'

# Row 17
$ws.Range("C17").Value = 86
Set-TextCell $ws 'D17' '9.0'
Set-TextCell $ws 'E17' 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F17").Value = 219
Set-TextCell $ws 'G17' '18834.00'

# Row 18
$ws.Range("C18").Value = 66
Set-TextCell $ws 'D18' '10.0'
Set-TextCell $ws 'E18' 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F18").Value = 303
Set-TextCell $ws 'G18' '19998.00'

# Row 19
Set-TextCell $ws 'A19' ''
$ws.Range("C19").Value = 70
Set-TextCell $ws 'D19' '11.0'
Set-TextCell $ws 'E19' 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F19").Value = 0
Set-TextCell $ws 'G19' '0.00'

# Row 20
Set-TextCell $ws 'A20' 'R. mtr.'
$ws.Range("C20").Value = 89
Set-TextCell $ws 'D20' '16'
Set-TextCell $ws 'E20' '20 mm'
$ws.Range("F20").Value = 40
Set-TextCell $ws 'G20' '3560.00'

# Row 21
$ws.Range("C21").Value = 37
Set-TextCell $ws 'D21' '17'
Set-TextCell $ws 'E21' '25 mm'
$ws.Range("F21").Value = 56
Set-TextCell $ws 'G21' '2072.00'

# Row 22
Set-TextCell $ws 'A22' ''
$ws.Range("C22").Value = 66
Set-TextCell $ws 'D22' '12.0'
Set-TextCell $ws 'E22' 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F22").Value = 0
Set-TextCell $ws 'G22' '0.00'

# Row 23
Set-TextCell $ws 'A23' 'Mtr.'
$ws.Range("C23").Value = 88
Set-TextCell $ws 'D23' '19'
Set-TextCell $ws 'E23' '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F23").Value = 81
Set-TextCell $ws 'G23' '7128.00'

# Row 24
$ws.Range("C24").Value = 7
Set-TextCell $ws 'D24' '20'
Set-TextCell $ws 'E24' '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F24").Value = 122
Set-TextCell $ws 'G24' '854.00'

# Row 25
$ws.Range("C25").Value = 57
Set-TextCell $ws 'G25' '326781.00'

# Row 26
$ws.Range("C26").Value = 70

# Row 27
Set-TextCell $ws 'A27' 'Mtr.'
$ws.Range("C27").Value = 90
Set-TextCell $ws 'D27' '23'
Set-TextCell $ws 'E27' '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F27").Value = 20
Set-TextCell $ws 'G27' '1800.00'

# Row 28
Set-TextCell $ws 'A28' ''
$ws.Range("C28").Value = 51
Set-TextCell $ws 'D28' '15.0'
Set-TextCell $ws 'E28' 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F28").Value = 0
Set-TextCell $ws 'G28' '0.00'

# Row 29
$ws.Range("C29").Value = 39

# Row 30
Set-TextCell $ws 'A30' ''
$ws.Range("C30").Value = 38
Set-TextCell $ws 'D30' '17.0'
Set-TextCell $ws 'E30' 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F30").Value = 0
Set-TextCell $ws 'G30' '0.00'

# Row 31
$ws.Range("C31").Value = 62
Set-TextCell $ws 'D31' '29'
Set-TextCell $ws 'E31' 'Single pole MCB   (With B/C curve tripping Characteristics)'

# Row 32
Set-TextCell $ws 'A32' 'Each'
$ws.Range("C32").Value = 13
Set-TextCell $ws 'D32' '30'
Set-TextCell $ws 'E32' ' 6 A to 32 A rating'
$ws.Range("F32").Value = 187
Set-TextCell $ws 'G32' '2431.00'

# Row 33
$ws.Range("C33").Value = 91

# Row 34
$ws.Range("C34").Value = 60
Set-TextCell $ws 'G34' '54000.00'

# Row 35
$ws.Range("C35").Value = 55

# Row 36
$ws.Range("C36").Value = 56

# Row 37
$ws.Range("C37").Value = 83
Set-TextCell $ws 'G37' '181272.00'

# Row 38
$ws.Range("C38").Value = 10

# Row 39
$ws.Range("C39").Value = 78

# Row 40
$ws.Range("C40").Value = 49

# Row 42
Set-TextCell $ws 'G42' '639218.00'
Set-TextCell $ws 'H42' '639218.00'

# Row 44
Set-TextCell $ws 'G44' '639218.00'
Set-TextCell $ws 'H44' '639218.00'
